$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 12.5045
$ws.Range("E6").Value = 11.8242
$ws.Range("D7").Value = -7.189299999999998
$ws.Range("B8").Value = 4.537100000000003
$ws.Range("E9").Value = 9.725099999999989
$ws.Range("B10").Value = 8.583900000000005
$ws.Range("E10").Value = 11.5666
$ws.Range("B12").Value = 5.808800000000002
$ws.Range("C13").Value = -12.8073
$ws.Range("B18").Value = 4.988900000000006
$ws.Range("D20").Value = -8.283199999999997
